$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Remove the "NO_LABEL" label cell (C3). The shared-string table will
# automatically drop the now-unused "NO_LABEL" entry and reindex the rest,
# which is the bulk of this diff.
$ws.Range("C3").Clear()

# Restore the previously-recorded selection on the survey sheet.
$ws.Range("C4").Select()
